# New board layout. Made things smaller, removed large trenches, and
# calculated USB differential impedance.
#
# The USB micro-B receptacle (J1, row 7) changed parts: the manufacturer,
# manufacturer part number, and LCSC part reference all need updating to
# reflect the new component.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D7").Value2 = "Amphenol ICC (FCI)"
$ws.Range("E7").Value2 = "10118193-0001LF"
$ws.Range("I7").Value2 = "LCSC Part: C132562"
